$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: drop the remix suffix from the title.
$ws.Range("A5").Value = "Bring Me To Life"

# Row 8: title/author were swapped in the source sheet - fix the mapping and
# drop the remix suffix from the title at the same time.
$ws.Range("A8").Value = " Flight"
$ws.Range("B8").Value = "Tristam & Braken"

# Add a hyperlink on the URL cell for row 2 (previously plain text).
$ws.Hyperlinks.Add($ws.Range("E2"), "https://youtu.be/yJg-Y5byMMw")
$ws.Range("E2").Style = "Hyperlink"

# Update the remembered selection to A9 (last data row).
$ws.Range("A9").Select()
